$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("E3", "E4", "E6", "E7", "C11", "C12", "C14", "D14", "C15", "B18", "C18", "B19", "C19", "B20", "C20", "B21", "C21")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = "-"
}
